# Apply "corrections to labs 1-3" edits to the Fick's Principle data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set explicit custom widths for columns A and B (previously default width).
# ColumnWidth is expressed in character units; the underlying file stores an
# offset value, so we compensate by subtracting that fixed offset (5/6) to
# land on the desired stored widths of 23 and ~16.57 characters.
$ws.Columns.Item(1).ColumnWidth = 22.166666666666668
$ws.Columns.Item(2).ColumnWidth = 15.736979166666666

# Correct the swapped Pulmonary Vein / Pulmonary Artery O2 values.
$ws.Range("B2").Value = 0.2
$ws.Range("B3").Value = 0.15
